# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New B:E values per row (F column - Win - stays unchanged); G is the row sum B+C+D+E
$data = @{
    2 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
    3 = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732)
    4 = @(3.230985683306322,  1.667794583268128,  0.8054896365839992, 0.496779210170732)
    5 = @(3.230985683306322,  1.667794583268128,  26.21740644021617,  8.660232485948974)
    6 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
    7 = @(3.230985683306322,  1.667794583268128,  0.8054896365839992, 0.496779210170732)
    8 = @(1.459612070389937,  1.667794583268128,  0.8054896365839992, 0.496779210170732)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = ($b + $c + $d + $e)
}
